$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Timestamp" header in column H, row 1
$ws.Range("H1").Value = "Timestamp"

# Update the active selection to H2, matching the post-edit state
$ws.Range("H2").Select()
